$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.4398919999999999
$ws.Range("I2").Value = 0.719974794695429
$ws.Range("J2").Value = 0.794097205716256
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.004751933329999999
$ws.Range("R2").Value = 0.02851159998
$ws.Range("S2").Value = 0.719974794695429
$ws.Range("T2").Value = 0.794097205716256

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1710905
$ws.Range("H3").Value = 0.342181
$ws.Range("I3").Value = 0.2800252053045709
$ws.Range("J3").Value = 0.2059027942837441
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 0.00184820512625
$ws.Range("R3").Value = 0.007392820505
$ws.Range("S3").Value = 0.2800252053045709
$ws.Range("T3").Value = 0.2059027942837441
